$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 (plain data row, default formatting) ---
$ws.Range("A13").Value = "Polychronidou P"
$ws.Range("B13").Value = "Kydros D"
$ws.Range("C13").Value = "lalala"
$ws.Range("D13").Value = 22222
$ws.Range("E13").Value = 2022
$ws.Range("F13").Value = "operations research"

# --- Row 14 ---
# A14 reuses the same "highlighted" cell format already used on B11/B12,
# so copy that formatting over before writing the new value.
$ws.Range("B11").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").Value = "Vrana V"

$ws.Range("B14").Value = "Kydros D"
$ws.Range("C14").Value = "lalala"
$ws.Range("D14").Value = 22222
$ws.Range("E14").Value = 2022
$ws.Range("F14").Value = "operations research"

# Move the selection to reflect where the user ended up after entering data
$ws.Range("B15").Select() | Out-Null
